$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HK_G_acc_LT"

$values = @(
    72.612612612612608,
    71.531531531531527,
    73.513513513513516,
    68.468468468468473,
    70.990990990990994,
    70.810810810810807,
    78.918918918918919,
    79.099099099099107,
    78.198198198198199,
    79.63963963963964,
    66.666666666666657,
    70.450450450450447,
    75.855855855855864,
    75.13513513513513,
    75.13513513513513,
    63.243243243243242,
    66.486486486486484,
    70.090090090090087,
    76.576576576576571,
    78.198198198198199,
    76.936936936936945,
    65.765765765765778,
    63.063063063063062,
    62.702702702702709,
    73.873873873873876,
    70.270270270270274,
    70.990990990990994,
    70.990990990990994,
    67.747747747747738,
    70.630630630630634,
    72.612612612612608,
    73.693693693693689,
    73.873873873873876,
    67.567567567567565,
    68.108108108108112,
    61.981981981981981,
    69.369369369369366,
    64.684684684684697,
    65.765765765765778,
    71.891891891891888,
    74.414414414414409,
    74.234234234234236,
    74.054054054054049,
    74.054054054054049,
    74.954954954954957,
    67.387387387387392,
    67.747747747747738,
    71.351351351351354
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

